$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook gained two new weekly records (rows 74-75) for
# "Vega Monumental Concepción" / Mandarina / Murcott, pushing the
# previously existing rows 74-84 down to rows 76-86.
$ws.Rows("74:75").Insert()

# --- Row 74 ---
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44491
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100102
$ws.Cells.Item(74, 8).Value = "Cítricos"
$ws.Cells.Item(74, 9).Value = 100102004
$ws.Cells.Item(74, 10).Value = "Mandarina"
$ws.Cells.Item(74, 11).Value = "Murcott"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 350
$ws.Cells.Item(74, 14).Value = 5500
$ws.Cells.Item(74, 15).Value = 6500
$ws.Cells.Item(74, 16).Value = 6071
$ws.Cells.Item(74, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(74, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(74, 19).Value = 337
$ws.Cells.Item(74, 20).Value = 18

# --- Row 75 ---
$ws.Cells.Item(75, 1).Value = 11
$ws.Cells.Item(75, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(75, 3).Value = "Bíobío"
$ws.Cells.Item(75, 4).Value = 44491
$ws.Cells.Item(75, 5).Value = 8
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100102
$ws.Cells.Item(75, 8).Value = "Cítricos"
$ws.Cells.Item(75, 9).Value = 100102004
$ws.Cells.Item(75, 10).Value = "Mandarina"
$ws.Cells.Item(75, 11).Value = "Murcott"
$ws.Cells.Item(75, 12).Value = "Segunda"
$ws.Cells.Item(75, 13).Value = 250
$ws.Cells.Item(75, 14).Value = 4000
$ws.Cells.Item(75, 15).Value = 4500
$ws.Cells.Item(75, 16).Value = 4300
$ws.Cells.Item(75, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(75, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(75, 19).Value = 239
$ws.Cells.Item(75, 20).Value = 18
